{"js": "const NEW_VALUES = [\"61+37=98\", \"49-44=5\", \"57-20=37\", \"21-4=17\", \"99-35=64\", \"37+2=39\", \"73+26=99\", \"0+36=36\", \"79+7=86\", \"71-56=15\", \"53-26=27\", \"95+2=97\", \"79+20=99\", \"55-4=51\", \"12+58=70\", \"88-54=34\", \"60+19=79\", \"12+20=32\", \"4+89=93\", \"38-26=12\", \"58+0=58\", \"82-28=54\", \"30-25=5\", \"34+32=66\", \"65-18=47\", \"71-16=55\", \"47+23=70\", \"7+22=29\", \"86-66=20\", \"91+7=98\", \"90-69=21\", \"79-70=9\", \"54+25=79\", \"90-81=9\", \"0+56=56\", \"99+0=99\", \"76-52=24\", \"80-56=24\", \"58+29=87\", \"69-30=39\", \"74-55=19\", \"94-80=14\", \"73+13=86\", \"88-47=41\", \"52-51=1\", \"28+0=28\", \"92-75=17\", \"5+6=11\", \"39-18=21\", \"99-36=63\", \"28+20=48\", \"35+18=53\", \"4+11=15\", \"44+2=46\", \"13+6=19\", \"70-55=15\", \"93-7=86\", \"39+29=68\", \"99-21=78\", \"8-0=8\", \"17+76=93\", \"34+22=56\", \"91-6=85\", \"17+64=81\", \"2+58=60\", \"88-5=83\", \"31+6=37\", \"69-58=11\", \"27+54=81\", \"96-72=24\", \"27+1=28\", \"90-88=2\", \"65+18=83\", \"0+3=3\", \"28+16=44\", \"49+29=78\", \"93-27=66\", \"86-33=53\", \"75-55=20\", \"11+6=17\", \"61-37=24\", \"58+40=98\", \"87-26=61\", \"61-40=21\", \"33-16=17\", \"2-1=1\", \"34+63=97\", \"5+73=78\", \"10-10=0\", \"42+40=82\", \"89-35=54\", \"98-3=95\", \"95-66=29\", \"97-45=52\", \"32+54=86\", \"66-21=45\", \"10+64=74\", \"18+35=53\", \"79-79=0\", \"12+4=16\"];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"No table found in document body\");\n}\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nconst colCount = table.values[0].length;\n\nfor (let i = 0; i < NEW_VALUES.length; i++) {\n  const row = Math.floor(i / colCount);\n  const col = i % colCount;\n  const cell = table.getCell(row, col);\n  cell.value = NEW_VALUES[i];\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$newValues = @(\"61+37=98\",\"49-44=5\",\"57-20=37\",\"21-4=17\",\"99-35=64\",\"37+2=39\",\"73+26=99\",\"0+36=36\",\"79+7=86\",\"71-56=15\",\"53-26=27\",\"95+2=97\",\"79+20=99\",\"55-4=51\",\"12+58=70\",\"88-54=34\",\"60+19=79\",\"12+20=32\",\"4+89=93\",\"38-26=12\",\"58+0=58\",\"82-28=54\",\"30-25=5\",\"34+32=66\",\"65-18=47\",\"71-16=55\",\"47+23=70\",\"7+22=29\",\"86-66=20\",\"91+7=98\",\"90-69=21\",\"79-70=9\",\"54+25=79\",\"90-81=9\",\"0+56=56\",\"99+0=99\",\"76-52=24\",\"80-56=24\",\"58+29=87\",\"69-30=39\",\"74-55=19\",\"94-80=14\",\"73+13=86\",\"88-47=41\",\"52-51=1\",\"28+0=28\",\"92-75=17\",\"5+6=11\",\"39-18=21\",\"99-36=63\",\"28+20=48\",\"35+18=53\",\"4+11=15\",\"44+2=46\",\"13+6=19\",\"70-55=15\",\"93-7=86\",\"39+29=68\",\"99-21=78\",\"8-0=8\",\"17+76=93\",\"34+22=56\",\"91-6=85\",\"17+64=81\",\"2+58=60\",\"88-5=83\",\"31+6=37\",\"69-58=11\",\"27+54=81\",\"96-72=24\",\"27+1=28\",\"90-88=2\",\"65+18=83\",\"0+3=3\",\"28+16=44\",\"49+29=78\",\"93-27=66\",\"86-33=53\",\"75-55=20\",\"11+6=17\",\"61-37=24\",\"58+40=98\",\"87-26=61\",\"61-40=21\",\"33-16=17\",\"2-1=1\",\"34+63=97\",\"5+73=78\",\"10-10=0\",\"42+40=82\",\"89-35=54\",\"98-3=95\",\"95-66=29\",\"97-45=52\",\"32+54=86\",\"66-21=45\",\"10+64=74\",\"18+35=53\",\"79-79=0\",\"12+4=16\")\n\n$table = $d.Tables.Item(1)\n$colCount = $table.Columns.Count\n\nfor ($i = 0; $i -lt $newValues.Count; $i++) {\n    $row = [math]::Floor($i / $colCount) + 1\n    $col = ($i % $colCount) + 1\n    $cell = $table.Cell($row, $col)\n    $cell.Range.Text = $newValues[$i]\n}\n"}
